$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H6").Value = "bit rate/s"
$ws.Range("I6").Value = "bit rate/s"
$ws.Range("H7").Value = "0xA53C000"
$ws.Range("I7").Value = 173260800
$ws.Range("I8").Formula = "=16*752*480*30"
